# Domain & UI 업데이트
# Insert a new "website: String" field paragraph into the Member domain
# class box on slide 2, immediately before the existing
# "interests: List<Hashtag>" line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# All the domain boxes on this slide live inside one big group shape.
$grp = $s.Shapes.Item(1)

# Find the "Member" class rectangle (shape Id 4) within the group.
$target = $null
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($sh.Id -eq 4) {
        $target = $sh
    }
}

$tr = $target.TextFrame.TextRange

# Locate the "interests: List<Hashtag>" paragraph and insert the new
# "website: String" paragraph directly before it, matching the
# surrounding run formatting (inherited automatically by InsertBefore).
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.TrimEnd("`r") -eq "interests: List<Hashtag>") {
        [void]$para.InsertBefore("website: String`r")
        break
    }
}
